$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the threshold values for alpha/beta/ratio ranges
$ws.Range("B2").Value = 5.6
$ws.Range("B3").Value = 5.7
$ws.Range("B4").Value = 0.7
$ws.Range("C4").Value = 1.3

# Remove the theta_threshold_range row entirely (row 5); the former
# pie_threshold_range row shifts up to become the new row 5.
$ws.Rows.Item(5).Delete()

# Update the (now shifted-up) pie_threshold_range row values
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 15

# Update the selection shown in the sheet view
[void]$ws.Range("B2:C3").Select()

# Set print/page setup properties
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
